$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: average of |S*|/n column (J)
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# New summary rows 14-17
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# Styling for B14:B17 - bold, size 12, vertical center aligned.
# Build the combined style on a scratch cell first (so only a single
# finished style ends up referenced in cellXfs), then paste just the
# formatting onto the target range and clean up the scratch cell.
$scratch = $ws.Range("AA1")
$scratch.Font.Bold = $true
$scratch.Font.Size = 12
$scratch.VerticalAlignment = -4108
$scratch.Copy()
$rng = $ws.Range("B14:B17")
$rng.PasteSpecial(-4122)
$scratch.Clear()

# Row heights to match the taller formatted rows
$ws.Rows.Item(14).RowHeight = 15.6
$ws.Rows.Item(15).RowHeight = 15.6
$ws.Rows.Item(16).RowHeight = 15.6
$ws.Rows.Item(17).RowHeight = 15.6

# Page setup
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Selection matching the saved view state
$ws.Range("A14:B17").Select()

$wb.Save()
